# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'44.413.16"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "'2.246.30"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'308.36"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "'94.78"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'35.05"
$ws.Range("E10").Value = "  +1.58%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "'7.22"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("E13").Value = "  +1.31%  "
$ws.Range("D14").Value = "'2.360.08"
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").Value = "'0.841"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "'13.72"
$ws.Range("E16").Value = "  +1.21%  "
$ws.Range("D17").Value = "'44.110.68"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'0.0₃0967"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'12.36"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D20").Value = "'6.42"
$ws.Range("E20").Value = "  +3.59%  "
$ws.Range("D21").Value = "'65.88"
$ws.Range("E21").Value = "  +1.85%  "
$ws.Range("D22").Value = "'3.20"
$ws.Range("E22").Value = "  +9.95%  "
$ws.Range("D23").Value = "'237.47"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "'2.02"
$ws.Range("E24").Value = "  +4.55%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'38.60"
$ws.Range("E26").Value = "  +7.03%  "
$ws.Range("D27").Value = "'2.23"
$ws.Range("E27").Value = "  +5.33%  "
$ws.Range("D28").Value = "'9.88"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").Value = "'5.96"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "'20.08"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "'153.91"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("D32").Value = "'0.0802"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "'3.12"
$ws.Range("E34").Value = "  -6.34%  "
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("D38").Value = "'3.51"
$ws.Range("E38").Value = "  +7.08%  "
$ws.Range("D39").Value = "'14.68"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "'3.83"
$ws.Range("E40").Value = "  +1.55%  "
$ws.Range("D41").Value = "'0.0305"
$ws.Range("E41").Value = "  +1.38%  "
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "'1.748.80"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("E44").Value = "  +3.59%  "
$ws.Range("D45").Value = "'80.63"
$ws.Range("E45").Value = "  -5.75%  "
$ws.Range("D46").Value = "'100.15"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").Value = "'70.98"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").Value = "'56.32"
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("E50").Value = "  +6.27%  "
$ws.Range("D51").Value = "'8.13"
$ws.Range("E51").Value = "  +0.35%  "
